# PLANILLA DE EVALUACION FINAL FASE 2 - fill in group names and rubric scores
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION1")

# Student names
$ws.Range("B4").Value = "CHRISTOPH BORNHARDT"
$ws.Range("B5").Value = "JOAN JARA"

# Grupal rubric (rows 13-19) - shared evaluation used for the whole team
$ws.Range("C13").Value = "Logrado"
$ws.Range("C14").Value = "Logro incipiente"
$ws.Range("C15").Value = "Logrado"
$ws.Range("C16").Value = "Logrado"
$ws.Range("C17").Value = "Logrado"
$ws.Range("C18").Value = "Logro incipiente"
$ws.Range("C19").Value = "Logro incipiente"

# Student 1 (CHRISTOPH BORNHARDT) individual rubric (rows 28-30)
$ws.Range("C28").Value = "Logrado"
$ws.Range("C29").Value = "Logrado"
$ws.Range("C30").Value = "Logrado"

# Student 2 (JOAN JARA) individual rubric (rows 40-42)
$ws.Range("C40").Value = "Logrado"
$ws.Range("C41").Value = "Logro incipiente"
$ws.Range("C42").Value = "Logrado"

# Update the selection / view to match the saved state
$ws.Range("L39").Select()
